# Applies the weekly update: inserts two new daily records at the top of the
# data block (rows 797 and 798), pushing all existing records down by two
# rows (797-831 -> 799-833). The dimension grows from A1:R831 to A1:R833.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 797; this shifts the old
# rows 797-831 down to 799-833 automatically.
$ws.Rows("797:798").Insert()

# New data for the two inserted rows (columns A..R).
$newRows = @(
    @(797, 3, "Femacal de La Calera", "Coquimbo", 45147, 5, 100112003, "Ajo", "Chino", "Primera", 68, 18000, 18500, 18221, "`$/caja 10 kilos", "China", 1822, 10, "Hortaliza"),
    @(798, 3, "Femacal de La Calera", "Coquimbo", 45147, 5, 100112003, "Ajo", "Chino", "Primera", 35, 21000, 21000, 21000, "`$/malla 10 kilos", "China", 2100, 10, "Hortaliza")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
    $ws.Cells.Item($rowNum, 9).Value = $r[9]
    $ws.Cells.Item($rowNum, 10).Value = $r[10]
    $ws.Cells.Item($rowNum, 11).Value = $r[11]
    $ws.Cells.Item($rowNum, 12).Value = $r[12]
    $ws.Cells.Item($rowNum, 13).Value = $r[13]
    $ws.Cells.Item($rowNum, 14).Value = $r[14]
    $ws.Cells.Item($rowNum, 15).Value = $r[15]
    $ws.Cells.Item($rowNum, 16).Value = $r[16]
    $ws.Cells.Item($rowNum, 17).Value = $r[17]
    $ws.Cells.Item($rowNum, 18).Value = $r[18]
}
